$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-order the product rows.
#
#    Rows 22-25 rotate up by one (row23's data moves to row22, row24's data
#    moves to row23, row25's data moves to row24, and row22's original data
#    wraps around into row25). Columns A:N hold the product data; column O
#    (timestamp) is handled separately below since every row gets the same
#    new value anyway.
#
#    A staging row (200) is used to hold the row-22 data while the rotation
#    happens, since the destination ranges overlap the sources. Each
#    destination is cleared immediately before the copy so that genuinely
#    blank source cells correctly blank out the destination (Copy() onto a
#    non-empty cell does not by itself clear a cell that has no counterpart
#    value in the source).
# ---------------------------------------------------------------------------

$ws.Range("A200:N200").Clear()
$ws.Range("A22:N22").Copy($ws.Range("A200:N200"))

$ws.Range("A22:N22").Clear()
$ws.Range("A23:N23").Copy($ws.Range("A22:N22"))

$ws.Range("A23:N23").Clear()
$ws.Range("A24:N24").Copy($ws.Range("A23:N23"))

$ws.Range("A24:N24").Clear()
$ws.Range("A25:N25").Copy($ws.Range("A24:N24"))

$ws.Range("A25:N25").Clear()
$ws.Range("A200:N200").Copy($ws.Range("A25:N25"))

$ws.Range("A200:N200").Clear()

# ---------------------------------------------------------------------------
# 2) Swap rows 51 and 52.
# ---------------------------------------------------------------------------

$ws.Range("A201:N201").Clear()
$ws.Range("A51:N51").Copy($ws.Range("A201:N201"))

$ws.Range("A51:N51").Clear()
$ws.Range("A52:N52").Copy($ws.Range("A51:N51"))

$ws.Range("A52:N52").Clear()
$ws.Range("A201:N201").Copy($ws.Range("A52:N52"))

$ws.Range("A201:N201").Clear()

# ---------------------------------------------------------------------------
# 3) Refresh the scrape timestamp (column O) for every data row.
# ---------------------------------------------------------------------------

$ws.Range("O2:O94").Value = "2022-08-02 20:57:25"
